$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed as a number by Excel;
# force Text number format first so the written cell keeps string semantics,
# matching the workbook's existing inlineStr/text-typed Price column.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D10", "D13", "D16", "D18", "D19", "D24", "D25", "D26", "D29", "D33", "D39", "D40", "D44", "D45", "D46", "D47", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value2 = "30.321.80"
$ws.Range("E2").Value2 = "  +1.48%  "

$ws.Range("D3").Value2 = "1.682.71"
$ws.Range("E3").Value2 = "  +3.42%  "

$ws.Range("D4").Value2 = "0.998"
$ws.Range("E4").Value2 = "  -0.33%  "

$ws.Range("D5").Value2 = "221.55"
$ws.Range("E5").Value2 = "  +3.53%  "

$ws.Range("D6").Value2 = "0.523"
$ws.Range("E6").Value2 = "  +0.70%  "

$ws.Range("D7").Value2 = "0.998"
$ws.Range("E7").Value2 = "  -0.24%  "

$ws.Range("D8").Value2 = "30.01"
$ws.Range("E8").Value2 = "  +1.38%  "

$ws.Range("E9").Value2 = "  +2.26%  "

$ws.Range("D10").Value2 = "0.0623"
$ws.Range("E10").Value2 = "  +1.69%  "

$ws.Range("E11").Value2 = "  -1.59%  "

$ws.Range("D12").Value2 = "1.922.24"
$ws.Range("E12").Value2 = "  +3.36%  "

$ws.Range("D13").Value2 = "10.73"
$ws.Range("E13").Value2 = "  +15.71%  "

$ws.Range("E14").Value2 = "  +8.59%  "

$ws.Range("D15").Value2 = "1.684.05"
$ws.Range("E15").Value2 = "  +3.44%  "

$ws.Range("D16").Value2 = "4.01"
$ws.Range("E16").Value2 = "  +3.60%  "

$ws.Range("D17").Value2 = "30.344.40"
$ws.Range("E17").Value2 = "  +1.52%  "

$ws.Range("D18").Value2 = "65.82"
$ws.Range("E18").Value2 = "  +1.42%  "

$ws.Range("D19").Value2 = "247.12"
$ws.Range("E19").Value2 = "  -0.38%  "

$ws.Range("D20").Value2 = "0.0₃0721"
$ws.Range("E20").Value2 = "  +2.40%  "

$ws.Range("E21").Value2 = "  -0.18%  "

$ws.Range("E22").Value2 = "  +3.79%  "

$ws.Range("E23").Value2 = "  +6.17%  "

$ws.Range("D24").Value2 = "2.20"
$ws.Range("E24").Value2 = "  +4.19%  "

$ws.Range("D25").Value2 = "158.78"
$ws.Range("E25").Value2 = "  -0.16%  "

$ws.Range("D26").Value2 = "15.88"
$ws.Range("E26").Value2 = "  +1.11%  "

$ws.Range("E27").Value2 = "  +0.34%  "

$ws.Range("E28").Value2 = "  +2.54%  "

$ws.Range("D29").Value2 = "0.998"
$ws.Range("E29").Value2 = "  -0.40%  "

$ws.Range("E30").Value2 = "  +2.34%  "

$ws.Range("E31").Value2 = "  +4.12%  "

$ws.Range("E32").Value2 = "  +1.28%  "

$ws.Range("D33").Value2 = "3.32"
$ws.Range("E33").Value2 = "  +3.56%  "

$ws.Range("D34").Value2 = "1.500.19"
$ws.Range("E34").Value2 = "  +5.06%  "

$ws.Range("E35").Value2 = "  +5.06%  "

$ws.Range("E36").Value2 = "  +0.03%  "

$ws.Range("E37").Value2 = "  +5.55%  "

$ws.Range("E38").Value2 = "  -4.46%  "

$ws.Range("D39").Value2 = "0.586"
$ws.Range("E39").Value2 = "  +5.99%  "

$ws.Range("D40").Value2 = "78.70"
$ws.Range("E40").Value2 = "  +11.09%  "

$ws.Range("E41").Value2 = "  +1.24%  "

$ws.Range("E42").Value2 = "  +2.89%  "

$ws.Range("E43").Value2 = "  +2.34%  "

$ws.Range("D44").Value2 = "2.00"
$ws.Range("E44").Value2 = "  +1.74%  "

$ws.Range("D45").Value2 = "0.998"
$ws.Range("E45").Value2 = "  -0.16%  "

$ws.Range("D46").Value2 = "1.00"
$ws.Range("E46").Value2 = "  -4.50%  "

$ws.Range("D47").Value2 = "51.82"
$ws.Range("E47").Value2 = "  -6.80%  "

$ws.Range("D48").Value2 = "1.813.68"
$ws.Range("E48").Value2 = "  +2.60%  "

$ws.Range("E49").Value2 = "  -0.30%  "

$ws.Range("D50").Value2 = "95.29"
$ws.Range("E50").Value2 = "  +6.30%  "

$ws.Range("D51").Value2 = "0.0₆0116"
$ws.Range("E51").Value2 = "  +8.90%  "
